$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 211, shifting existing rows 211..292 down to 212..293
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with the new record.
# Columns A,B,C,E,F,G,H,I,J,K are constant across this dataset (copy the template values).
$ws.Cells.Item(211, 1).Value2 = 4
$ws.Cells.Item(211, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(211, 3).Value2 = "Los Lagos"
$ws.Cells.Item(211, 4).Value2 = 44489
$ws.Cells.Item(211, 5).Value2 = 10
$ws.Cells.Item(211, 6).Value2 = "Fruta"
$ws.Cells.Item(211, 7).Value2 = 100102
$ws.Cells.Item(211, 8).Value2 = "Cítricos"
$ws.Cells.Item(211, 9).Value2 = 100102003
$ws.Cells.Item(211, 10).Value2 = "Limón"
$ws.Cells.Item(211, 11).Value2 = "Sin especificar"
$ws.Cells.Item(211, 12).Value2 = "1a amarillo"
$ws.Cells.Item(211, 13).Value2 = 200
$ws.Cells.Item(211, 14).Value2 = 9000
$ws.Cells.Item(211, 15).Value2 = 9000
$ws.Cells.Item(211, 16).Value2 = 9000
$ws.Cells.Item(211, 17).Value2 = "$/malla 18 kilos"
$ws.Cells.Item(211, 18).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(211, 19).Value2 = 500
$ws.Cells.Item(211, 20).Value2 = 18

# Make sure the date cell keeps the same date style as other date cells (s="2"),
# which the row-insert operation already carried down from row 210.
$ws.Cells.Item(211, 4).NumberFormat = $ws.Cells.Item(210, 4).NumberFormat
